# Add procedural-comment test text to the two Kaiwá cells for "two" (row 5)
# and move the active selection of the bottom-right frozen pane to G11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H5").Value = "<mokõi>(dos){Guasch1962:616} (PCP: Test procedural comment lands in variants)"
$ws.Range("J5").Value = "[mõˈkõj̃] (PCP: Procedural comment before real comment) (dos)"

$null = $ws.Range("G11").Select()
